# Update workbook "al 08/04/2020": append the new monthly inflation row
# (March 2020, date serial 43896 -> B88 = 21.2) right after the last
# existing data row (row 87), mirroring its formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 88

# Copy the previous data row (A87:B87 -> A88:B88) so the new row inherits
# the same date / number styles used throughout the table, then overwrite
# the copied values with the new month's figures.
$ws.Range("A87:B87").Copy($ws.Range("A88:B88")) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 43896
$ws.Cells.Item($newRow, 2).Value = 21.2

# Move the active selection to F82, matching the updated saved selection.
$ws.Range("F82").Select() | Out-Null
